$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$chData = New-Object 'object[,]' 30,6
$chData[0,0] = -5.570007801055908
$chData[0,1] = -35.83580017089844
$chData[0,2] = 22.41282081604004
$chData[0,3] = -0.03355996282946583
$chData[0,4] = -0.1584969647691545
$chData[0,5] = -0.0141987274892259
$chData[1,0] = -2.182258129119873
$chData[1,1] = -7.31383228302002
$chData[1,2] = 33.1801872253418
$chData[1,3] = 2.232942344609228
$chData[1,4] = -2.220868759447247
$chData[1,5] = -1.187766118819664
$chData[2,0] = -1.57672917842865
$chData[2,1] = -10.6974925994873
$chData[2,2] = 15.88597106933594
$chData[2,3] = 1.551477467791645
$chData[2,4] = 0.9860338054088031
$chData[2,5] = -2.580682175500044
$chData[3,0] = 9.968832015991213
$chData[3,1] = -22.49884223937988
$chData[3,2] = 30.88003540039062
$chData[3,3] = -2.007183966429317
$chData[3,4] = 3.929720215175458
$chData[3,5] = 0.6425367645595088
$chData[4,0] = -6.721622467041016
$chData[4,1] = -1.676235437393188
$chData[4,2] = -11.5780611038208
$chData[4,3] = -5.864261772321608
$chData[4,4] = 5.580999001212737
$chData[4,5] = 2.13303363693427
$chData[5,0] = -32.24742889404297
$chData[5,1] = -20.76686477661133
$chData[5,2] = 7.779660701751709
$chData[5,3] = -1.424133648783631
$chData[5,4] = -2.927099547771023
$chData[5,5] = 1.214256319940471
$chData[6,0] = -7.299997329711914
$chData[6,1] = -29.29754447937012
$chData[6,2] = 42.08726501464844
$chData[6,3] = 7.334702552475461
$chData[6,4] = -7.529888727650148
$chData[6,5] = -2.842065592730211
$chData[7,0] = -13.7584810256958
$chData[7,1] = -9.336126327514648
$chData[7,2] = 20.6390266418457
$chData[7,3] = 4.214317862661773
$chData[7,4] = -3.208229954198307
$chData[7,5] = -4.081179956471704
$chData[8,0] = 24.05830383300781
$chData[8,1] = -37.75975799560547
$chData[8,2] = 69.02017974853516
$chData[8,3] = -2.042688489330464
$chData[8,4] = 3.253492230213918
$chData[8,5] = -1.052057455785544
$chData[9,0] = -65.68665313720703
$chData[9,1] = -20.87689781188965
$chData[9,2] = -14.51115989685059
$chData[9,3] = -8.89284183371872
$chData[9,4] = 6.183227533139049
$chData[9,5] = 8.303981626996203
$chData[10,0] = -3.672505617141724
$chData[10,1] = -5.974431037902832
$chData[10,2] = 3.619078159332275
$chData[10,3] = 1.101257697395411
$chData[10,4] = -8.100705602894388
$chData[10,5] = 4.845630360686418
$chData[11,0] = -15.47988033294678
$chData[11,1] = -32.08433532714844
$chData[11,2] = 16.76184463500977
$chData[11,3] = 6.093993388347799
$chData[11,4] = -8.757804864682157
$chData[11,5] = 1.959031731445559
$chData[12,0] = 9.497495651245115
$chData[12,1] = -6.459360122680664
$chData[12,2] = 12.91139125823975
$chData[12,3] = 5.412965199962155
$chData[12,4] = -0.5916063237634659
$chData[12,5] = -1.592202567165649
$chData[13,0] = -4.22635555267334
$chData[13,1] = -2.0812087059021
$chData[13,2] = 10.98652267456055
$chData[13,3] = 1.925117330521559
$chData[13,4] = 2.407113150780237
$chData[13,5] = -2.267917759485144
$chData[14,0] = 12.04420471191406
$chData[14,1] = -23.34278678894043
$chData[14,2] = 30.86811256408692
$chData[14,3] = -2.498231235498205
$chData[14,4] = 2.901872266153379
$chData[14,5] = 1.981700023694026
$chData[15,0] = -40.70948791503906
$chData[15,1] = 10.1422061920166
$chData[15,2] = -4.144575119018555
$chData[15,3] = -4.308581977157101
$chData[15,4] = 6.365314957517947
$chData[15,5] = -0.1142475042283233
$chData[16,0] = -17.44337844848633
$chData[16,1] = -33.51113891601562
$chData[16,2] = -9.5015230178833
$chData[16,3] = -1.724872384752597
$chData[16,4] = 9.335920878819055
$chData[16,5] = -6.735600778034694
$chData[17,0] = -12.57210731506348
$chData[17,1] = -66.97329711914062
$chData[17,2] = 49.08029174804688
$chData[17,3] = 2.584916405055822
$chData[17,4] = -0.09003212140948058
$chData[17,5] = -6.952868482340953
$chData[18,0] = -9.933971405029297
$chData[18,1] = -2.512709140777588
$chData[18,2] = 6.143205642700195
$chData[18,3] = 3.838851762854538
$chData[18,4] = -1.53637689981414
$chData[18,5] = -4.963172649004442
$chData[19,0] = 21.00619125366211
$chData[19,1] = -11.01493453979492
$chData[19,2] = 19.27869033813477
$chData[19,3] = 1.802320775223111
$chData[19,4] = 0.8488230512749766
$chData[19,5] = -0.6512930319176782
$chData[20,0] = -6.04071855545044
$chData[20,1] = -12.84956741333008
$chData[20,2] = -4.132137298583984
$chData[20,3] = -4.243776159812123
$chData[20,4] = 0.1656510622604825
$chData[20,5] = 6.680505474161589
$chData[21,0] = 34.18217468261719
$chData[21,1] = 1.69527006149292
$chData[21,2] = -2.47593355178833
$chData[21,3] = -3.268159332482785
$chData[21,4] = 6.903478548393394
$chData[21,5] = 2.484958121495333
$chData[22,0] = -27.50531387329102
$chData[22,1] = -56.38927841186523
$chData[22,2] = -12.15304946899414
$chData[22,3] = 1.061862719725376
$chData[22,4] = 5.545756017014141
$chData[22,5] = -2.003967119299857
$chData[23,0] = 41.86726379394531
$chData[23,1] = -60.87873458862305
$chData[23,2] = -19.13644981384277
$chData[23,3] = 5.571206705910729
$chData[23,4] = -0.27861049012367
$chData[23,5] = 0.5144404683794082
$chData[24,0] = -7.41383171081543
$chData[24,1] = 0.0885686874389648
$chData[24,2] = 9.766034126281738
$chData[24,3] = 7.170254665872369
$chData[24,4] = -1.238433612429565
$chData[24,5] = -0.1954094016033365
$chData[25,0] = -51.14746856689453
$chData[25,1] = -84.02630615234375
$chData[25,2] = 86.47578430175781
$chData[25,3] = 1.018974213866787
$chData[25,4] = 1.231977397610628
$chData[25,5] = -0.4505931842401376
$chData[26,0] = -39.8954963684082
$chData[26,1] = -2.818742275238037
$chData[26,2] = -19.23029327392578
$chData[26,3] = -9.593666403930271
$chData[26,4] = -3.09781322716174
$chData[26,5] = -1.724484086777074
$chData[27,0] = 1.352512359619141
$chData[27,1] = -3.199073314666748
$chData[27,2] = 9.475702285766602
$chData[27,3] = -4.399059822840986
$chData[27,4] = -10.81391397322185
$chData[27,5] = -1.643179933847163
$chData[28,0] = -6.012904644012451
$chData[28,1] = -19.74169921875
$chData[28,2] = 2.084371328353882
$chData[28,3] = 4.520313150394069
$chData[28,4] = -8.092394739944785
$chData[28,5] = 0.05753792109697597
$chData[29,0] = 10.77586936950684
$chData[29,1] = -30.89251708984375
$chData[29,2] = 17.9653205871582
$chData[29,3] = 5.967537736300274
$chData[29,4] = -5.414584435291156
$chData[29,5] = -1.474100663795217

$ws.Range("C2:H31").Value = $chData

# New rows 22-31 need timestamp (A) and label (B) values too
$abData = New-Object 'object[,]' 10,2
$abData[0,0] = 2000
$abData[0,1] = "walkingToRunning"
$abData[1,0] = 2100
$abData[1,1] = "walkingToRunning"
$abData[2,0] = 2200
$abData[2,1] = "walkingToRunning"
$abData[3,0] = 2300
$abData[3,1] = "walkingToRunning"
$abData[4,0] = 2400
$abData[4,1] = "walkingToRunning"
$abData[5,0] = 2500
$abData[5,1] = "walkingToRunning"
$abData[6,0] = 2600
$abData[6,1] = "walkingToRunning"
$abData[7,0] = 2700
$abData[7,1] = "walkingToRunning"
$abData[8,0] = 2800
$abData[8,1] = "walkingToRunning"
$abData[9,0] = 2900
$abData[9,1] = "walkingToRunning"
$ws.Range("A22:B31").Value = $abData
